# edit.ps1 - apply the changes described by the target diff:
#   1. Slide 6's table changes its table style GUID
#      from {56EC123C-1495-42E1-A18B-A48711C9F3DD}
#      to   {BD549A98-20E9-41D7-8164-FB3B16F083CB}
#   2. The deck's theme colour scheme (ppt/theme/theme1.xml, the theme
#      actually used by the slide master / slides) is switched from the
#      "Integral" palette to the stock "Office Theme" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 -------------------------------------------
$s   = $p.Slides.Item(6)
$sh  = $s.Shapes.Item(2)
$tbl = $sh.Table
$tbl.ApplyStyle("{BD549A98-20E9-41D7-8164-FB3B16F083CB}")

# --- 2. Theme colours: Integral -> Office Theme --------------------------
$m  = $p.SlideMaster
$cs = $m.ColorScheme

$cs.Item(1).RGB  = 0          # dk1      000000
$cs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$cs.Item(3).RGB  = 6968388    # dk2      44546A
$cs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$cs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$cs.Item(6).RGB  = 3243501    # accent2  ED7D31
$cs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$cs.Item(8).RGB  = 49407      # accent4  FFC000
$cs.Item(9).RGB  = 12874308   # accent5  4472C4
$cs.Item(10).RGB = 4697456    # accent6  70AD47
$cs.Item(11).RGB = 12673797   # hlink    0563C1
$cs.Item(12).RGB = 7491477    # folHlink 954F72
